# Update "想去人数" (want-to-go count) figures for several events that
# appear both on their category sheet ("展览") and on the aggregate
# "全部类型" sheet, matching freshly regenerated site data.

$wb = $excel.ActiveWorkbook

# -- Sheet "展览" (Exhibition) --
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 3640
$wsExpo.Range("F5").Value = 2230
$wsExpo.Range("F11").Value = 1340
$wsExpo.Range("F13").Value = 1994

# -- Sheet "全部类型" (All types) --
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 3640
$wsAll.Range("F5").Value = 2230
$wsAll.Range("F14").Value = 1340
$wsAll.Range("F16").Value = 1994
